$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 287; this shifts existing rows 287:395 down to 288:396
$ws.Rows("287:287").Insert()

# Populate the newly inserted row 287 with its data
$ws.Range("A287").Value = 3
$ws.Range("B287").Value = "Femacal de La Calera"
$ws.Range("C287").Value = "Coquimbo"
$ws.Range("D287").Value = 44795
$ws.Range("E287").Value = 5
$ws.Range("F287").Value = 100112009
$ws.Range("G287").Value = "Acelga"
$ws.Range("H287").Value = "Sin especificar"
$ws.Range("I287").Value = "Primera"
$ws.Range("J287").Value = 250
$ws.Range("K287").Value = 3300
$ws.Range("L287").Value = 3400
$ws.Range("M287").Value = 3348
$ws.Range("N287").Value = "$/docena de atados (6 kilos)"
$ws.Range("O287").Value = "Provincia de Quillota"
$ws.Range("P287").Value = 558
$ws.Range("Q287").Value = 6
$ws.Range("R287").Value = "Hortaliza"
